$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 640 (shifts existing rows 640.. down to 642..)
$ws.Rows("640:641").Insert()

# --- New row 640 ---
$ws.Range("A640").Value = 6
$ws.Range("B640").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C640").Value = "Metropolitana"
$ws.Range("D640").Value = 45124
$ws.Range("E640").Value = 13
$ws.Range("F640").Value = 100112039
$ws.Range("G640").Value = "Ciboulette"
$ws.Range("H640").Value = "Sin especificar"
$ws.Range("I640").Value = "Primera"
$ws.Range("J640").Value = 200
$ws.Range("K640").Value = 2000
$ws.Range("L640").Value = 2000
$ws.Range("M640").Value = 2000
$ws.Range("N640").Value = "`$/docena de atados"
$ws.Range("O640").Value = "Región Metropolitana"
$ws.Range("P640").Value = 667
$ws.Range("Q640").Value = 3
$ws.Range("R640").Value = "Hortaliza"

# --- New row 641 ---
$ws.Range("A641").Value = 6
$ws.Range("B641").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C641").Value = "Metropolitana"
$ws.Range("D641").Value = 45124
$ws.Range("E641").Value = 13
$ws.Range("F641").Value = 100112039
$ws.Range("G641").Value = "Ciboulette"
$ws.Range("H641").Value = "Sin especificar"
$ws.Range("I641").Value = "Segunda"
$ws.Range("J641").Value = 300
$ws.Range("K641").Value = 1500
$ws.Range("L641").Value = 1500
$ws.Range("M641").Value = 1500
$ws.Range("N641").Value = "`$/docena de atados"
$ws.Range("O641").Value = "Región Metropolitana"
$ws.Range("P641").Value = 500
$ws.Range("Q641").Value = 3
$ws.Range("R641").Value = "Hortaliza"

# Apply the same date number format style as the other D-column cells to the two new D cells
$ws.Range("D640:D641").NumberFormat = $ws.Range("D642").NumberFormat
